# API: Gameweeks import (#25)
# Adds two new columns ("Show Statistics Continuously" and "Gameweek")
# to the end of the "Challenges" sheet, with sample data in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Challenges")

# New header cells
$ws.Range("S1").Value = "Show Statistics Continuously"
$ws.Range("T1").Value = "Gameweek"

# New sample data cells
# NOTE: leading apostrophe forces Excel to store this as literal text
# ("true") rather than auto-converting it to a Boolean TRUE; resetting
# the style afterwards drops the "quote prefix" formatting flag that
# the apostrophe entry would otherwise leave behind.
$ws.Range("S2").Value = "'true"
$ws.Range("S2").Style = "Normal"
$ws.Range("T2").Value = 1
